# Auto-generated edit script: updates cryptos list values per commit
# "Updated cryptos list on Tue Feb 21 17:23:47 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.516.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.51%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.99%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3900"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.56%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3926"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.002"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.28%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.397"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08623"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.270"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001316"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.699"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.673.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07059"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.045"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.70%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.504.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.378"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.23%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.737"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.57%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.811"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -14.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "147.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.38%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.263"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.523"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.52%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.859.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.32%  "

# Row 34
$ws.Range("E34").Value = "  -5.47%  "

# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.943"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.41%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03014"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2800"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "

# Row 38
$ws.Range("E38").Value = "  -3.78%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09422"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.533"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7883"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.43%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.62%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7088"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.540"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.170"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.81%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08566"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.317"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.28%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.32%  "
